$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.070.15"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "3.052.12"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.046.22"
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.17%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.451"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "3.540.50"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Value = "63.123.89"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").Value = "3.050.70"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.669"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "462.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0811"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").Value = "3.060.54"
$ws.Range("E40").Value = "  -5.34%  "
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("E51").Value = "  +0.86%  "
